$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells I1 ("I0") and J1 ("IF") ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (style) of the existing header cell H1 onto the new
# header cells so they match the bold/centered/bordered header style (s="1").
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data values for columns I (I0) and J (IF), rows 2-59 ---
$iValues = @(8,6,7,7,5,7,8,10,5,8,8,4,8,8,7,6,6,7,6,5,7,9,7,6,8,8,7,8,6,8,7,7,8,6,7,8,4,7,6,8,7,9,7,6,6,8,12,8,7,6,7,7,6,6,6,6,5,5)
$jValues = @(8,6,7,7,5,7,8,10,6,8,8,5,8,8,7,7,6,7,6,6,7,9,7,6,8,8,8,8,6,8,7,7,8,6,7,8,4,7,7,8,7,9,7,7,6,8,12,8,7,6,7,7,6,6,6,6,5,5)

for ($r = 2; $r -le 59; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
